$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# CaseField tab: insert a new "CategoryID" column before the existing
# "FieldTypeParameter" column (column H).
# ---------------------------------------------------------------------------
$wsCaseField = $wb.Worksheets.Item("CaseField")
$wsCaseField.Columns.Item(8).Insert()
$wsCaseField.Cells.Item(2, 8).Value = "`t`nA non-mandatory field. This field will be used to indicate a category for a document field type or Collection of Document field ONLY. `nMust match to a valid CategoryID defined in the Categories tab for the given CaseTypeID."
$wsCaseField.Cells.Item(3, 8).Value = "CategoryID"

# ---------------------------------------------------------------------------
# ComplexTypes tab: insert a new "CategoryID" column before the existing
# "FieldTypeParameter" column (column F).
# ---------------------------------------------------------------------------
$wsComplexTypes = $wb.Worksheets.Item("ComplexTypes")
$wsComplexTypes.Columns.Item(6).Insert()
$wsComplexTypes.Cells.Item(2, 6).Value = "A non-mandatory field. This field will be used to indicate a category for a document field type or Collection of Document field ONLY. `nMust match to a valid CategoryID defined in the Categories tab."
$wsComplexTypes.Cells.Item(3, 6).Value = "CategoryID"

# ---------------------------------------------------------------------------
# Restore a selection/active-tab state similar to the authored edit: leave
# ComplexTypes as the active sheet (selection on the new CategoryID cell).
# ---------------------------------------------------------------------------
$wsCaseField.Activate()
$wsCaseField.Range("Q2").Select()

$wsComplexTypes.Activate()
$wsComplexTypes.Range("F3").Select()
